# Add a new "Save" column (H) to the sheet, mirroring the other stat
# columns (B..G) that already live in row 1 / rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: text "Save", formatted like the rest of the header row.
$ws.Cells.Item(1, 8).Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for H2:H37 (the "Save" flag for each game), in row order.
$saveValues = @(0,0,0,0,0,1,0,0,0,1,0,0,1,0,1,0,0,1,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,0,0,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
